# Update Inscritos (E), Pagos (F) and Inscrições homologadas (H) figures
# for the affected course/campus rows in the "Inscricoes" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 10; E = 419; F = 209; H = 209 },
    @{ Row = 11; E = 281; F = 154; H = 154 },
    @{ Row = 12; E = 408 },
    @{ Row = 13; E = 108 },
    @{ Row = 14; E = 104 },
    @{ Row = 15; E = 136 },
    @{ Row = 17; E = 79 },
    @{ Row = 22; E = 147 },
    @{ Row = 23; E = 175 },
    @{ Row = 24; E = 180 },
    @{ Row = 25; E = 219; F = 100; H = 100 },
    @{ Row = 26; E = 123 },
    @{ Row = 27; E = 279 },
    @{ Row = 28; E = 167 },
    @{ Row = 29; E = 146 },
    @{ Row = 30; E = 181; F = 100; H = 100 },
    @{ Row = 34; E = 181; F = 111; H = 111 },
    @{ Row = 35; E = 120 },
    @{ Row = 36; E = 56; F = 34; H = 34 },
    @{ Row = 41; E = 336; F = 157; H = 157 },
    @{ Row = 42; E = 307 },
    @{ Row = 43; E = 103; F = 51; H = 51 },
    @{ Row = 44; E = 260 },
    @{ Row = 45; E = 119 },
    @{ Row = 46; E = 266; F = 144; H = 144 },
    @{ Row = 47; E = 377 },
    @{ Row = 48; E = 180; F = 72; H = 72 },
    @{ Row = 49; E = 254; F = 104; H = 104 },
    @{ Row = 52; E = 22; F = 10; H = 10 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Range("E$r").Value = $u.E
    if ($u.ContainsKey("F")) { $ws.Range("F$r").Value = $u.F }
    if ($u.ContainsKey("H")) { $ws.Range("H$r").Value = $u.H }
}
